$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.175.71'
$ws.Range("E2").Value = '  +2.57%  '

$ws.Range("D3").Value = '2.347.86'
$ws.Range("E3").Value = '  +6.37%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").Value = '''109.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.57%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '''312.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.40%  '

$ws.Range("D7").Value = '''0.643'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.23%  '

$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("E9").Value = '  +7.10%  '

$ws.Range("D10").Value = '''42.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.62%  '

$ws.Range("D11").Value = '''0.0938'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.35%  '

$ws.Range("E12").Value = '  +2.18%  '

$ws.Range("E13").Value = '  +12.05%  '

$ws.Range("E14").Value = '  +2.53%  '

$ws.Range("D15").Value = '''16.27'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.49%  '

$ws.Range("D16").Value = '2.702.89'
$ws.Range("E16").Value = '  +6.51%  '

$ws.Range("D17").Value = '2.349.73'
$ws.Range("E17").Value = '  +5.26%  '

$ws.Range("D18").Value = '43.124.40'
$ws.Range("E18").Value = '  +2.64%  '

$ws.Range("E19").Value = '  +4.27%  '

$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("D21").Value = '''75.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.52%  '

$ws.Range("E22").Value = '  +14.34%  '

$ws.Range("D23").Value = '''3.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.47%  '

$ws.Range("D24").Value = '''253.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.75%  '

$ws.Range("D25").Value = '''9.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.37%  '

$ws.Range("D26").Value = '''12.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.95%  '

$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").Value = '''39.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.91%  '

$ws.Range("E29").Value = '  +0.69%  '

$ws.Range("D30").Value = '''22.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.28%  '

$ws.Range("D31").Value = '''174.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.20%  '

$ws.Range("E32").Value = '  -0.62%  '

$ws.Range("D33").Value = '''0.0929'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.48%  '

$ws.Range("D34").Value = '''5.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.86%  '

$ws.Range("E35").Value = '  +6.72%  '

$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("E37").Value = '  +5.15%  '

$ws.Range("D38").Value = '''4.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.38%  '

$ws.Range("E39").Value = '  +2.64%  '

$ws.Range("D40").Value = '''2.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.27%  '

$ws.Range("D41").Value = '''72.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.27%  '

$ws.Range("E42").Value = '  +15.26%  '

$ws.Range("E43").Value = '  +2.50%  '

$ws.Range("D44").Value = '''12.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.90%  '

$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("E46").Value = '  +4.68%  '

$ws.Range("E47").Value = '  +11.30%  '

$ws.Range("D48").Value = '''110.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.51%  '

$ws.Range("E49").Value = '  -0.72%  '

$ws.Range("E50").Value = '  +3.27%  '

$ws.Range("D51").Value = '''70.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.54%  '

